# Updated symbol list on Sun Jan 22 04:29:14 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$value) {
    # Prefix with an apostrophe so Excel stores the exact text
    # (no auto-conversion to a number/percentage), then strip the
    # leftover quote-prefix style so the cell keeps its original styling.
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("E2") "-0.82%"

Set-TextValue $ws.Range("D3") "38.09"
Set-TextValue $ws.Range("E3") "8.66%"

Set-TextValue $ws.Range("D4") "4.977"
Set-TextValue $ws.Range("E4") "-3.33%"

Set-TextValue $ws.Range("D5") "0.07732"
Set-TextValue $ws.Range("E5") "-0.46%"

Set-TextValue $ws.Range("D6") "2.185"
Set-TextValue $ws.Range("E6") "-7.24%"

Set-TextValue $ws.Range("D7") "8.004"
Set-TextValue $ws.Range("E7") "-0.53%"

Set-TextValue $ws.Range("E8") "1.27%"

Set-TextValue $ws.Range("D9") "0.9134"
Set-TextValue $ws.Range("E9") "-1.70%"

Set-TextValue $ws.Range("D10") "0.09274"
Set-TextValue $ws.Range("E10") "-7.67%"

Set-TextValue $ws.Range("E11") "-0.13%"

Set-TextValue $ws.Range("D12") "0.08406"
Set-TextValue $ws.Range("E12") "-2.61%"

Set-TextValue $ws.Range("D13") "0.03544"
Set-TextValue $ws.Range("E13") "6.87%"

Set-TextValue $ws.Range("E14") "0.17%"

Set-TextValue $ws.Range("D15") "0.001479"
Set-TextValue $ws.Range("E15") "-1.27%"

Set-TextValue $ws.Range("D16") "0.005727"
Set-TextValue $ws.Range("E16") "-1.29%"

Set-TextValue $ws.Range("D17") "3.472"
Set-TextValue $ws.Range("E17") "0.27%"

Set-TextValue $ws.Range("E18") "2.64%"

Set-TextValue $ws.Range("E19") "3.09%"

Set-TextValue $ws.Range("D20") "0.1316"
Set-TextValue $ws.Range("E20") "-1.31%"

Set-TextValue $ws.Range("D21") "4.541"
Set-TextValue $ws.Range("E21") "5.64%"

Set-TextValue $ws.Range("D22") "0.2230"
Set-TextValue $ws.Range("E22") "-3.14%"

Set-TextValue $ws.Range("D23") "0.04652"
Set-TextValue $ws.Range("E23") "1.94%"

Set-TextValue $ws.Range("D24") "0.001228"
Set-TextValue $ws.Range("E24") "1.34%"

Set-TextValue $ws.Range("D25") "0.004437"
Set-TextValue $ws.Range("E25") "1.54%"

Set-TextValue $ws.Range("D26") "0.0001299"
Set-TextValue $ws.Range("E26") "0.01%"

Set-TextValue $ws.Range("D27") "0.0004747"
Set-TextValue $ws.Range("E27") "39.84%"

Set-TextValue $ws.Range("D39") "0.01749"
Set-TextValue $ws.Range("E39") "-2.32%"

Set-TextValue $ws.Range("D40") "0.04686"
Set-TextValue $ws.Range("E40") "-2.50%"

Set-TextValue $ws.Range("D41") "0.007868"
Set-TextValue $ws.Range("E41") "1.20%"

Set-TextValue $ws.Range("D42") "0.1387"
Set-TextValue $ws.Range("E42") "-1.85%"

Set-TextValue $ws.Range("D43") "0.007649"
Set-TextValue $ws.Range("E43") "4.28%"

Set-TextValue $ws.Range("D44") "0.002289"
Set-TextValue $ws.Range("E44") "7.02%"

Set-TextValue $ws.Range("D45") "0.009954"
Set-TextValue $ws.Range("E45") "5.30%"

Set-TextValue $ws.Range("D46") "0.00006056"
Set-TextValue $ws.Range("E46") "-0.86%"

Set-TextValue $ws.Range("D47") "0.00000000750"
Set-TextValue $ws.Range("E47") "0.01%"

Set-TextValue $ws.Range("D48") "8.691"
Set-TextValue $ws.Range("E48") "183.51%"

Set-TextValue $ws.Range("E49") "35.01%"

Set-TextValue $ws.Range("D50") "0.00002099"
Set-TextValue $ws.Range("E50") "0.01%"

Set-TextValue $ws.Range("D51") "0.0001999"
Set-TextValue $ws.Range("E51") "0.01%"

